# Split the trailing run of the "Topic overview" paragraph into three runs,
# inserting the word "in" so the sentence reads:
#   "...I am also interested in variance in those variables..."
# This mirrors how Word stores an in-place mid-sentence word insertion as
# separate <w:r> runs around the newly typed text.

$d = $word.ActiveDocument
$content = $d.Content

$anchor = "I am also interested variance in those variables between SNOTEL sites at similar elevations."
$found = $content.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence to edit."
}

$target = $content.Paragraphs(1).Range

$expected = "`tTopic overview: The dataset that I intend to work with is the SNOTEL dataset. This is an automated system of snowpack and climate sensors which collect data on snow water equivalent (SWE), precipitation, and air temperature. Some sensors within the SNOTEL network also collect other variables such as snow depth, soil moisture, and wind speed. This dataset is important because it provides consistent data on snowpack, especially in areas that are difficult or costly to access for field surveys. SNOTEL is also important because it can be used to project annual water supply, predict floods, and it also has been used in climate studies. I intend to use this dataset to investigate variance in SWE, precipitation, and temperature at SNOTEL sites across different EPA Level III Ecoregions. I am also interested variance in those variables between SNOTEL sites at similar elevations.`r"
if ($target.Text -ne $expected) {
    throw "Paragraph text did not match the expected original content; aborting to avoid corrupting the document."
}

$newParagraphXml = '<w:p w14:paraId="5861D388" w14:textId="4A5CB42A" w:rsidR="00E319D8" w:rsidRDefault="006C7A2D" w:rsidP="00AB13B5"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:tab/></w:r><w:r w:rsidR="00A17913"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Topic </w:t></w:r><w:r w:rsidR="00A0425B"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>overview:</w:t></w:r><w:r w:rsidR="00A0425B"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00A803A5"><w:t>The dataset that I intend to work with is the SNOTEL dataset. This is an automated system of snowpack and climate sensors which collect data</w:t></w:r><w:r w:rsidR="00AB13B5"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="005B034E"><w:t>on snow water equivalent (SWE), precipitation, and air temperature. Some sensors within the SNOTEL network also collect other variables</w:t></w:r><w:r w:rsidR="003441A9"><w:t xml:space="preserve"> such as snow depth, soil moisture, and wind speed</w:t></w:r><w:r w:rsidR="005B034E"><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="003441A9"><w:t xml:space="preserve">This dataset is important because it provides </w:t></w:r><w:r w:rsidR="0060372D"><w:t>consistent data on snowpack, especially in areas that are difficult or costly to access for field surveys.</w:t></w:r><w:r w:rsidR="00E319D8"><w:t xml:space="preserve"> SNOTEL is also important because it can be used to project annual water supply, predict floods, and </w:t></w:r><w:r w:rsidR="00A17913"><w:t xml:space="preserve">it </w:t></w:r><w:r w:rsidR="00E319D8"><w:t>also has been used in climate studies.</w:t></w:r><w:r w:rsidR="00A17913"><w:t xml:space="preserve"> I intend to use this dataset to investigate variance in SWE, precipitation, and temperature at SNOTEL sites across different EPA Level III Ecoregions. I am also interested</w:t></w:r><w:r w:rsidR="00A17913"><w:t xml:space="preserve"> in</w:t></w:r><w:r w:rsidR="00A17913"><w:t xml:space="preserve"> variance in those variables between SNOTEL sites at similar elevations.</w:t></w:r></w:p>'
$target.InsertXML($newParagraphXml)
